# Update "想去人数" (want-to-go count) values in column F across sheets,
# reflecting a newly generated data snapshot (gh-pages output update).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2801
$ws1.Range("F3").Value = 1130
$ws1.Range("F4").Value = 20463
$ws1.Range("F6").Value = 2563
$ws1.Range("F8").Value = 614
$ws1.Range("F9").Value = 481
$ws1.Range("F14").Value = 396
$ws1.Range("F17").Value = 177
$ws1.Range("F18").Value = 238
$ws1.Range("F19").Value = 22
$ws1.Range("F20").Value = 24
$ws1.Range("F21").Value = 113

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F12").Value = 95

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6077
$ws3.Range("F5").Value = 1384
$ws3.Range("F6").Value = 39

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6077
$ws4.Range("F5").Value = 1384
$ws4.Range("F6").Value = 2801
$ws4.Range("F7").Value = 1130
$ws4.Range("F8").Value = 20464
$ws4.Range("F14").Value = 2564
$ws4.Range("F17").Value = 39
$ws4.Range("F18").Value = 614
$ws4.Range("F19").Value = 481
$ws4.Range("F27").Value = 396
$ws4.Range("F32").Value = 95
$ws4.Range("F33").Value = 177
$ws4.Range("F35").Value = 238
$ws4.Range("F39").Value = 22
$ws4.Range("F42").Value = 24
$ws4.Range("F48").Value = 113

$wb.Save()
